$d = $word.ActiveDocument

# 1. Update the "Last updated" date field result text.
$d.Content.Find.Execute("April 30, 2015", $true, $false, $false, $false, $false, $true, 1, $false, "April 25, 2020", 2)

# 2. Remove the stray "_GoBack" bookmark that sits right after the first
#    "Manual" run (Word drops/repositions this automatically as edits are
#    made elsewhere in the document).
$d.Bookmarks("_GoBack").Delete()

# 3. Replace "_rvtApp.Create.NewPlane" with "Plane.CreateByNormalAndOrigin"
#    on the "pGeomPlaneH" line. Re-assert the (identical) Courier New
#    formatting explicitly on the freshly typed text so it stays its own
#    run, matching how Word records an in-place retype, then drop Word's
#    "_GoBack" bookmark at the point where the edit finished (exactly
#    what Word itself does after the most recent edit).
$r = $d.Content
$r.Find.Execute("_rvtApp.Create.NewPlane")
$r.Text = "Plane.CreateByNormalAndOrigin"
$r.Font.Name = "Courier New"
$r.Font.NameAscii = "Courier New"
$r.Font.NameBi = "Courier New"
$r.Font.Size = 10
$r.NoProofing = 1
$ins = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $ins)

# 4. Bump the Revit SDK path year referenced in the sample App.config from
#    2013 to 2019.
$r2 = $d.Content
$r2.Find.Execute("Revit SDK 201")
$digit = $d.Range($r2.End, $r2.End + 1)
$digit.Text = "9"
